{"js": "// Xbox GDK Samples Readme \u2014 update references to \"Game Core\" APIs to the\n// November GDK release naming (\"WINAPI_FAMILY_GAMES\" / \"Microsoft GDK\"),\n// matching the commit \"Update Xbox GDK Samples to November GDK release.\"\n//\n// Each of the five occurrences of the phrase \"Game Core \" in the document\n// body is handled individually (search on a longer, unambiguous phrase so\n// we never touch the wrong paragraph), because each needs a different\n// replacement:\n//   1. \"...use of the Game Core XGameSave APIs...\"            -> delete \"Game Core \"\n//   2. \"...related to the new Game Core API surface area...\"  -> \"WINAPI_FAMILY_GAMES \"\n//   3. \"...working with the Game Core APIs. Most of...\"        -> \"Microsoft GDK \"\n//   4. \"RAII class wrappers around the Game Core handle...\"    -> delete \"Game Core \"\n//   5. \"...new task system, Game Core APIs, use synchronous...\" -> \"Microsoft GDK \"\n\nasync function replaceOnce(context, searchPhrase, replacement) {\n  const results = context.document.body.search(searchPhrase, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Search phrase not found: \" + searchPhrase);\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"This sample demonstrates the use of the Game Core XGameSave APIs\"\nawait replaceOnce(\n  context,\n  \"use of the Game Core \",\n  \"use of the \"\n);\n\n// 2) \"...a variety of other techniques related to the new Game Core API surface area...\"\nawait replaceOnce(\n  context,\n  \"related to the new Game Core API surface area\",\n  \"related to the new WINAPI_FAMILY_GAMES API surface area\"\n);\n\n// 3) \"...assist in working with the Game Core APIs. Most of these files...\"\nawait replaceOnce(\n  context,\n  \"working with the Game Core APIs\",\n  \"working with the Microsoft GDK APIs\"\n);\n\n// 4) \"RAII class wrappers around the Game Core handle types...\"\nawait replaceOnce(\n  context,\n  \"RAII class wrappers around the Game Core handle types\",\n  \"RAII class wrappers around the handle types\"\n);\n\n// 5) \"...rewritten to utilize new task system, Game Core APIs, use synchronous (blocking)...\"\nawait replaceOnce(\n  context,\n  \"new task system, Game Core APIs, use synchronous\",\n  \"new task system, Microsoft GDK APIs, use synchronous\"\n);\n", "ps1": "# Xbox GDK Samples Readme \u2014 update references to \"Game Core\" APIs to the\n# November GDK release naming (\"WINAPI_FAMILY_GAMES\" / \"Microsoft GDK\"),\n# matching the commit \"Update Xbox GDK Samples to November GDK release.\"\n#\n# Each of the five occurrences of the phrase \"Game Core \" in the document\n# body is handled individually (searching on a longer, unambiguous phrase\n# so we never touch the wrong paragraph), because each needs a different\n# replacement:\n#   1. \"...use of the Game Core XGameSave APIs...\"             -> delete \"Game Core \"\n#   2. \"...related to the new Game Core API surface area...\"   -> \"WINAPI_FAMILY_GAMES \"\n#   3. \"...working with the Game Core APIs. Most of...\"         -> \"Microsoft GDK \"\n#   4. \"RAII class wrappers around the Game Core handle...\"     -> delete \"Game Core \"\n#   5. \"...new task system, Game Core APIs, use synchronous...\" -> \"Microsoft GDK \"\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText($FindText, $ReplaceText) {\n    $rng = $d.Content\n    # wdFindContinue = 1, wdReplaceOne = 2\n    $found = $rng.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $FindText\"\n    }\n}\n\nReplace-DocText \"use of the Game Core \" \"use of the \"\n\nReplace-DocText \"related to the new Game Core API surface area\" \"related to the new WINAPI_FAMILY_GAMES API surface area\"\n\nReplace-DocText \"working with the Game Core APIs\" \"working with the Microsoft GDK APIs\"\n\nReplace-DocText \"RAII class wrappers around the Game Core handle types\" \"RAII class wrappers around the handle types\"\n\nReplace-DocText \"new task system, Game Core APIs, use synchronous\" \"new task system, Microsoft GDK APIs, use synchronous\"\n"}
